# Apply the change described in the commit:
# "added 'other' instead of NaN in the ghgrp to include all unit type"
#
# This adds a new "Other" row to the "unit types" sheet (mirroring the
# existing generic unit types), removes the stray leftover value in G14,
# and makes the "unit types" sheet the active/selected tab.

$wb = $excel.ActiveWorkbook

$wsTypes = $wb.Worksheets.Item("unit types")

# Remove the stray leftover value that lived far outside the real table.
$wsTypes.Range("G14").ClearContents()

# Add the new "Other" unit type row, matching the other "generic" rows.
$wsTypes.Range("A7").Value = "Other"
$wsTypes.Range("B7").Value = "generic"
$wsTypes.Range("D7").Value = "Y"

# Make "unit types" the active sheet/tab, with B7 selected.
$wsTypes.Activate()
$wsTypes.Range("B7").Select()
